# isyonetim.xlsx — add per-district tracking sheets, a couple of new
# entries, and tidy up a stray empty row.

$wb = $excel.ActiveWorkbook

$kayitlar = $wb.Worksheets.Item("Kayitlar")
$personel = $wb.Worksheets.Item("Personel")
$isler    = $wb.Worksheets.Item("İşler")
$birim    = $wb.Worksheets.Item("Birim")

# ------------------------------------------------------------------
# 1) Add six new district sheets after "Birim", each seeded with the
#    same header row as "Kayitlar" (Kayıt No / Tarih / Birim / Dosya
#    Sayısı / Parsel Sayısı / İş / Personeller).
# ------------------------------------------------------------------
$kayitlar.Range("A1:G1").Copy()

$districtNames = @("Merkez İlçe", "Anamur", "Silifke", "Erdemli", "Mut ", "Tarsus")
$prev = $birim
$newSheets = @{}
foreach ($name in $districtNames) {
    $sheet = $wb.Worksheets.Add($null, $prev)
    $sheet.Name = $name
    $sheet.Range("A1:G1").PasteSpecial()
    $newSheets[$name] = $sheet
    $prev = $sheet
}

# ------------------------------------------------------------------
# 2) "Birim" gains a row for the new "TARSUS" district.
# ------------------------------------------------------------------
$birim.Range("A6").Value = "TARSUS"

# ------------------------------------------------------------------
# 3) "Kayitlar" picks up a stray quote-mark entry in E16.
# ------------------------------------------------------------------
$kayitlar.Range("E16").Value = '"'

# ------------------------------------------------------------------
# 4) "İşler" loses a trailing blank formatted row (row 17).
# ------------------------------------------------------------------
[void]$isler.Rows("17:17").Delete()

# ------------------------------------------------------------------
# 5) Leave the UI state roughly where the author left it.
# ------------------------------------------------------------------
[void]$newSheets["Anamur"].Range("R15").Select()
[void]$newSheets["Silifke"].Range("P19").Select()
[void]$newSheets["Tarsus"].Range("G13").Select()
[void]$birim.Range("E18").Select()
[void]$isler.Range("D14").Select()

$kayitlar.Activate()
[void]$kayitlar.Range("O16").Select()
